$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update A2 and A4 to the new image filename "blog_4.jpg"
$ws.Range("A2").Value = "blog_4.jpg"
$ws.Range("A4").Value = "blog_4.jpg"

# Move the active selection to A4
$ws.Range("A4").Select()
